$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row for the "Single, 2-Hole 90deg Bracket" part above the
# existing "Single, 2-Hole Plate" row (was row 22, becomes row 23). ---
$ws.Rows.Item(22).Insert()

# Re-apply the thin-box border that every other data row uses, since a
# freshly inserted row comes in borderless.
$ws.Range("A22:E22").Borders.LineStyle = 1
$ws.Range("A22:E22").Borders.Weight = 2
$ws.Range("A22:E22").Borders.ColorIndex = 1

# Fill in the new row's data.
$ws.Range("A22").Value = "Single, 2-Hole 90" + [char]0x00B0 + " Bracket"
$ws.Range("B22").Value = "McMaster - Carr"
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = "5537T181"
$ws.Range("E22").Value = 31.04
$ws.Range("E22").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# The "Fasteners (4 pack)" row (previously row 23, now shifted to row 24)
# had its quantity updated from 28 to 34.
$ws.Range("C24").Value = 34

# Update the view: scroll / select like the saved workbook did.
$ws.Range("A36").Select()

# Touch the page setup so a pageSetup element is written for the sheet.
$ws.PageSetup.Orientation = 1

Write-Host "Bill of materials updated"
